$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data, copying the formatting of the row above it
# (same shaded style used by the other "page1-*" rows)
$ws.Range("A13:B13").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("A14").Value = "page1-children_seen"
$ws.Range("B14").Value = "nb_children_seen"

# Match the recorded selection from the saved workbook
$ws.Range("B16").Select()
